# "Update countries & provincias Spain"
# The "Pais" sheet lists COVID-19 stats per country (col A = country,
# B..H = Casos totales/Nuevos casos/Casos activos/Recuperados/
# Casos criticos/Muertes hoy/Muertes). This refreshes the snapshot to a
# later pull (14:50 -> 16:07) and updates the per-country figures that
# moved, including a handful of rows where the country at that rank
# changed (new data shuffled "Kenia"/"El Salvador"/"Gabon",
# "Mauritania"/"Republica de Africa Central", "Libano"/"Eslovenia",
# "Dominica"/"Fiyi" and "Islas Turcas y Caicos"/"Santa Sede").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 16:07"

# Plain data refresh (country/rank unchanged)
$ws.Range("B7").Value = 400566
$ws.Range("C7").Value = 4754
$ws.Range("D7").Value = 216676
$ws.Range("E7").Value = 170860
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = 13030

$ws.Range("B8").Value = 303110
$ws.Range("C8").Value = 1295
$ws.Range("G8").Value = 128
$ws.Range("H8").Value = 42589

$ws.Range("B37").Value = 38841
$ws.Range("C37").Value = 377
$ws.Range("D37").Value = 24906
$ws.Range("E37").Value = 12407
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 1528

$ws.Range("E50").Value = 5570
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 59

$ws.Range("B59").Value = 12803
$ws.Range("C59").Value = 94
$ws.Range("E59").Value = 1032
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 260

$ws.Range("B69").Value = 8733
$ws.Range("C69").Value = 7
$ws.Range("E69").Value = 351

$ws.Range("B79").Value = 5399
$ws.Range("C79").Value = 61
$ws.Range("D79").Value = 3894
$ws.Range("E79").Value = 1453
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 52

# Rows 84-86: Kenia / El Salvador / Gabon re-sorted
$ws.Range("A84").Value = "Kenia"
$ws.Range("B84").Value = 4478
$ws.Range("C84").Value = 104
$ws.Range("D84").Value = 1586
$ws.Range("E84").Value = 2773
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 119

$ws.Range("A85").Value = "El Salvador"
$ws.Range("B85").Value = 4475
$ws.Range("C85").Value = 146
$ws.Range("D85").Value = 2449
$ws.Range("E85").Value = 1933
$ws.Range("G85").Value = 7
$ws.Range("H85").Value = 93

$ws.Range("A86").Value = "Gabon"
$ws.Range("B86").Value = 4428
$ws.Range("D86").Value = 1750
$ws.Range("E86").Value = 2644
$ws.Range("H86").Value = 34

# Rows 97-98: Mauritania / Republica de Africa Central swapped
$ws.Range("A97").Value = "Republica de Africa Central"
$ws.Range("B97").Value = 2686
$ws.Range("C97").Value = 81
$ws.Range("D97").Value = 420
$ws.Range("E97").Value = 2247
$ws.Range("H97").Value = 19

$ws.Range("A98").Value = "Mauritania"
$ws.Range("B98").Value = 2621
$ws.Range("D98").Value = 653
$ws.Range("E98").Value = 1866
$ws.Range("H98").Value = 102

# Plain data refresh (country/rank unchanged)
$ws.Range("B100").Value = 2309
$ws.Range("C100").Value = 4
$ws.Range("D100").Value = 2071
$ws.Range("E100").Value = 153

$ws.Range("D106").Value = 1472
$ws.Range("E106").Value = 467

$ws.Range("B111").Value = 1822
$ws.Range("C111").Value = 3
$ws.Range("E111").Value = 11

# Rows 116-117: Eslovenia / Libano swapped
$ws.Range("A116").Value = "Libano"
$ws.Range("B116").Value = 1536
$ws.Range("C116").Value = 26
$ws.Range("D116").Value = 1006
$ws.Range("E116").Value = 498
$ws.Range("H116").Value = 32

$ws.Range("A117").Value = "Eslovenia"
$ws.Range("B117").Value = 1519
$ws.Range("C117").Value = 6
$ws.Range("D117").Value = 1359
$ws.Range("E117").Value = 51
$ws.Range("H117").Value = 109

# Plain data refresh (country/rank unchanged)
$ws.Range("B130").Value = 901
$ws.Range("C130").Value = 1
$ws.Range("E130").Value = 38

$ws.Range("B142").Value = 688
$ws.Range("C142").Value = 20
$ws.Range("E142").Value = 507

$ws.Range("B161").Value = 287
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 196
$ws.Range("E161").Value = 85

# Rows 202-203: Fiyi / Dominica swapped (no numeric changes)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Rows 208-209: Santa Sede / Islas Turcas y Caicos swapped
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
